$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (their new values look numeric)
$textCells = @("D5", "D6", "D10", "D11", "D12", "D16", "D18", "D20", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D33", "D38", "D39", "D42", "D44", "D45", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.919.82'
$ws.Range("E2").Value = '  +4.94%  '
$ws.Range("D3").Value = '2.268.07'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("D5").Value = '302.62'
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("D6").Value = '92.34'
$ws.Range("E6").Value = '  +6.31%  '
$ws.Range("E7").Value = '  +3.30%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +3.78%  '
$ws.Range("D10").Value = '54.54'
$ws.Range("E10").Value = '  +8.50%  '
$ws.Range("D11").Value = '32.26'
$ws.Range("E11").Value = '  +6.02%  '
$ws.Range("D12").Value = '0.0801'
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("E14").Value = '  +3.94%  '
$ws.Range("D15").Value = '2.619.07'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").Value = '14.16'
$ws.Range("E16").Value = '  +3.01%  '
$ws.Range("D17").Value = '2.262.40'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("D18").Value = '0.759'
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").Value = '41.853.33'
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").Value = '12.10'
$ws.Range("E20").Value = '  +8.71%  '
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("E22").Value = '  +3.53%  '
$ws.Range("D23").Value = '67.11'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '241.62'
$ws.Range("E24").Value = '  +1.76%  '
$ws.Range("E25").Value = '  +4.33%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '1.91'
$ws.Range("E27").Value = '  +4.01%  '
$ws.Range("D28").Value = '24.00'
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  +4.35%  '
$ws.Range("E30").Value = '  -4.98%  '
$ws.Range("D31").Value = '159.34'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").Value = '33.90'
$ws.Range("E32").Value = '  +6.50%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +4.21%  '
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("E37").Value = '  +2.48%  '
$ws.Range("D38").Value = '0.105'
$ws.Range("E38").Value = '  +5.52%  '
$ws.Range("D39").Value = '16.60'
$ws.Range("E39").Value = '  +9.32%  '
$ws.Range("E40").Value = '  +3.45%  '
$ws.Range("E41").Value = '  +4.76%  '
$ws.Range("D42").Value = '3.93'
$ws.Range("E42").Value = '  +6.18%  '
$ws.Range("D43").Value = '2.070.33'
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("D44").Value = '19.64'
$ws.Range("E44").Value = '  +9.60%  '
$ws.Range("D45").Value = '0.0280'
$ws.Range("E45").Value = '  +2.98%  '
$ws.Range("D46").Value = '10.14'
$ws.Range("E46").Value = '  +3.57%  '
$ws.Range("E47").Value = '  +6.67%  '
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("E49").Value = '  +3.71%  '
$ws.Range("E50").Value = '  +3.69%  '
$ws.Range("D51").Value = '51.84'
$ws.Range("E51").Value = '  +5.76%  '
